# RA2 Ch1 - updated roadmap
# as per feedback from call on 03/01/20

$p = $ppt.ActivePresentation

# 1. Update the roadmap SmartArt wording on slide 4 ("Diagram 1"):
#    "All Infrastructure Profiles included"
#      -> "Mapping of specification to RM Infrastructure Profiles included"
$roadmapSlide = $p.Slides.Item(4)
$diagramShape = $roadmapSlide.Shapes.Item(1)
$smartArt = $diagramShape.SmartArt
$oldText = "All Infrastructure Profiles included"
$newText = "Mapping of specification to RM Infrastructure Profiles included"
for ($i = 1; $i -le $smartArt.AllNodes.Count; $i++) {
    $node = $smartArt.AllNodes.Item($i)
    if ($node.TextFrame2.TextRange.Text -eq $oldText) {
        $node.TextFrame2.TextRange.Text = $newText
    }
}

# 2. Refresh the fixed "date updated" footer field (slide master + every
#    slide layout) from 24/12/2019 to 08/01/2020, matching a
#    Header&Footer "Apply to All" with a fixed date.
$newDate = "08/01/2020"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
